$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values look numeric/date-like to Excel's auto-detection,
# so force them to be treated as plain text (matching the source inlineStr cells)
# by setting NumberFormat to Text before assignment, then restore the default
# "Normal" style so no stray formatting is introduced.
$dCells = @('D2','D3','D4','D5','D6','D7','D8','D9','D10','D11','D13','D15','D16','D17','D18','D19','D20','D21','D22','D23','D25','D26','D27','D28','D29','D30','D31','D32','D33','D34','D35','D36','D37','D38','D39','D40','D42','D43','D44','D45','D46','D48','D49','D50','D51')
foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range('D2').Value = '30.410.69'
$ws.Range('D3').Value = '2.102.24'
$ws.Range('D4').Value = '1.003'
$ws.Range('D5').Value = '334.28'
$ws.Range('D6').Value = '1.002'
$ws.Range('D7').Value = '0.5221'
$ws.Range('D8').Value = '0.4555'
$ws.Range('D9').Value = '54.45'
$ws.Range('D10').Value = '0.08894'
$ws.Range('D11').Value = '1.179'
$ws.Range('D13').Value = '2.091.26'
$ws.Range('D15').Value = '8.026'
$ws.Range('D16').Value = '97.13'
$ws.Range('D17').Value = '0.00001145'
$ws.Range('D18').Value = '1.003'
$ws.Range('D19').Value = '0.06621'
$ws.Range('D20').Value = '19.18'
$ws.Range('D21').Value = '1.002'
$ws.Range('D22').Value = '6.297'
$ws.Range('D23').Value = '30.467.38'
$ws.Range('D25').Value = '2.358'
$ws.Range('D26').Value = '2.331.78'
$ws.Range('D27').Value = '22.19'
$ws.Range('D28').Value = '162.58'
$ws.Range('D29').Value = '2.516'
$ws.Range('D30').Value = '133.02'
$ws.Range('D31').Value = '1.207'
$ws.Range('D32').Value = '0.1067'
$ws.Range('D33').Value = '1.651'
$ws.Range('D34').Value = '6.391'
$ws.Range('D35').Value = '3.933'
$ws.Range('D36').Value = '10.40'
$ws.Range('D37').Value = '5.860'
$ws.Range('D38').Value = '0.02571'
$ws.Range('D39').Value = '0.06838'
$ws.Range('D40').Value = '0.2313'
$ws.Range('D42').Value = '0.6867'
$ws.Range('D43').Value = '1.248'
$ws.Range('D44').Value = '2.320'
$ws.Range('D45').Value = '13.98'
$ws.Range('D46').Value = '0.6352'
$ws.Range('D48').Value = '1.246'
$ws.Range('D49').Value = '0.00000000344'
$ws.Range('D50').Value = '83.15'
$ws.Range('D51').Value = '1.201'

foreach ($addr in $dCells) { $ws.Range($addr).Style = "Normal" }

# Volume(1h) column (E) values are already unambiguous text (percent strings
# with surrounding spaces), so they can be assigned directly.
$ws.Range('E2').Value = '  -0.27%  '
$ws.Range('E3').Value = '  -0.20%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('E5').Value = '  +1.34%  '
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('E7').Value = '  -0.85%  '
$ws.Range('E8').Value = '  +3.71%  '
$ws.Range('E9').Value = '  +14.38%  '
$ws.Range('E10').Value = '  +0.23%  '
$ws.Range('E12').Value = '  -2.16%  '
$ws.Range('E13').Value = '  -0.67%  '
$ws.Range('E14').Value = '  +0.87%  '
$ws.Range('E15').Value = '  +3.25%  '
$ws.Range('E16').Value = '  +0.72%  '
$ws.Range('E17').Value = '  +1.20%  '
$ws.Range('E19').Value = '  -0.30%  '
$ws.Range('E20').Value = '  +0.80%  '
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('E23').Value = '  -0.27%  '
$ws.Range('E24').Value = '  +0.20%  '
$ws.Range('E25').Value = '  +0.20%  '
$ws.Range('E26').Value = '  -0.84%  '
$ws.Range('E27').Value = '  -1.15%  '
$ws.Range('E29').Value = '  -3.52%  '
$ws.Range('E30').Value = '  +0.18%  '
$ws.Range('E31').Value = '  -0.41%  '
$ws.Range('E32').Value = '  -0.65%  '
$ws.Range('E33').Value = '  -1.67%  '
$ws.Range('E34').Value = '  +2.71%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('E36').Value = '  +2.23%  '
$ws.Range('E37').Value = '  +6.48%  '
$ws.Range('E38').Value = '  -0.47%  '
$ws.Range('E39').Value = '  +1.99%  '
$ws.Range('E40').Value = '  +1.35%  '
$ws.Range('E41').Value = '  -0.72%  '
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('E43').Value = '  -1.43%  '
$ws.Range('E44').Value = '  +4.71%  '
$ws.Range('E45').Value = '  -0.75%  '
$ws.Range('E46').Value = '  -0.61%  '
$ws.Range('E47').Value = '  +0.56%  '
$ws.Range('E48').Value = '  -0.55%  '
$ws.Range('E49').Value = '  +17.15%  '
$ws.Range('E50').Value = '  +0.97%  '
$ws.Range('E51').Value = '  -0.87%  '
